$wb = $excel.ActiveWorkbook

# Overview sheet: row for 530290e6-5491-4005-a69b-99d51fd2293c.md now reports
# "Handed back" status for both zh-cn and de-de columns (was "Ready for handoff").
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: same row's Status flips to Handed back, the handback datetime
# advances, and the previous "not latest" error is cleared now that the
# handback report is in sync.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-06 06:56:57"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: mirror of the same update.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-06 06:57:13"
$wsDeDe.Range("P3").Value = ""

# The Error Detail column no longer holds long text, so Excel's column
# autofit shrinks it back down from its fixed width of 40.
$wsZhCn.Columns.Item(16).AutoFit() | Out-Null
$wsDeDe.Columns.Item(16).AutoFit() | Out-Null
